$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C1").Value = "tech_reason"

$ws.Range("C2").Value = "The candidate lacks the required skills in MongoDB and NodeJS, which are crucial for the job. The projects demonstrate good proficiency in ReactJS and JavaScript, but the absence of MongoDB and NodeJS skills impacts the score."
$ws.Range("E2").Value = "The candidate has demonstrated strong adaptability, problem-solving skills, and a willingness to learn in the interview answers. Additionally, the candidate's interest in AI and the desire to contribute to cutting-edge technologies align well with the company's goals, but the preference to not work alone may need to be addressed."

$ws.Range("C3").Value = "The candidate lacks the required skills of MongoDB and Web Development, which are essential for the role. Although they have experience with some relevant technologies, the absence of these key skills has impacted the score."
$ws.Range("D3").Value = 3
$ws.Range("E3").Value = "The candidate demonstrates strong adaptability, teamwork, and problem-solving skills. Their eagerness to immerse in a new culture and work environment, along with a clear career plan, shows a positive attitude. The expressed interest in learning from the Japanese work culture and language also demonstrates flexibility and adaptability, making them a suitable candidate for the role."

$ws.Range("C4").Value = "The candidate lacks direct experience with MERN stack (MongoDB, ExpressJS, ReactJS, NodeJS) which are primary requirements for the job. However, the candidate's projects showcase strong skills in AI/ML, Python, and ReactJS, which could be beneficial in a tech-driven environment."
$ws.Range("D4").Value = 3
$ws.Range("E4").Value = "The candidate shows a strong interest in Japanese work culture and willingness to adapt to a new environment. The candidate also demonstrates good teamwork and problem-solving skills. However, there is room for improvement in time management and adaptability, which are important for working in a foreign country."

$wb.Save()
